# Updated cryptos list on Wed Mar  8 06:44:24 UTC 2023 with GitHub Actions
#
# Note: several "Price" (column D) values look like plain decimal numbers
# (e.g. "286.76"); Excel auto-converts those to numeric cells unless they
# are entered as text. We prefix those with a leading apostrophe (quote
# prefix) so they stay text, matching the source data which stores every
# Price cell as a string. Values that already contain two dots (e.g.
# "22.015.87") can't be parsed as a number anyway, so no prefix is needed
# there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "22.015.87"
$ws.Range("E2").Value = "  -1.94%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.554.10"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4 (TetherUSD) - only volume changes
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 (USDC) - only volume changes
$ws.Range("E5").Value = "  +0.01%  "

# Row 6 (BNB)
$ws.Range("D6").Value = "'286.76"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.3762"
$ws.Range("E7").Value = "  +1.11%  "

# Row 8 (Cardano)
$ws.Range("D8").Value = "'0.3237"
$ws.Range("E8").Value = "  -2.53%  "

# Row 9 / Row 10 swap: OKB <-> Polygon (identity swap with updated price/volume)
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.127"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'41.29"
$ws.Range("E10").Value = "  -13.15%  "

# Row 11 (Dogecoin)
$ws.Range("D11").Value = "'0.07298"
$ws.Range("E11").Value = "  -3.02%  "

# Row 12 (BinanceUSD)
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13 (Solana)
$ws.Range("D13").Value = "'19.61"
$ws.Range("E13").Value = "  -5.64%  "

# Row 14 (Polkadot)
$ws.Range("D14").Value = "'5.698"
$ws.Range("E14").Value = "  -4.00%  "

# Row 15 (Chainlink)
$ws.Range("D15").Value = "'6.844"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16 (WrappedEther) - only price changes
$ws.Range("D16").Value = "1.556.45"

# Row 17 (ShibaInu)
$ws.Range("D17").Value = "'0.00001081"
$ws.Range("E17").Value = "  -3.37%  "

# Row 18 (TRON)
$ws.Range("D18").Value = "'0.06641"
$ws.Range("E18").Value = "  -1.30%  "

# Row 19 (Litecoin)
$ws.Range("D19").Value = "'85.02"
$ws.Range("E19").Value = "  -3.73%  "

# Row 20 (Uniswap)
$ws.Range("D20").Value = "'6.440"
$ws.Range("E20").Value = "  +0.64%  "

# Row 21 (Dai)
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.02%  "

# Row 22 (Avalanche)
$ws.Range("D22").Value = "'15.96"
$ws.Range("E22").Value = "  -3.38%  "

# Row 23 (Cosmos)
$ws.Range("D23").Value = "'11.56"
$ws.Range("E23").Value = "  -3.70%  "

# Row 24 (WrappedBTC)
$ws.Range("D24").Value = "22.028.12"
$ws.Range("E24").Value = "  -1.81%  "

# Row 25 (Toncoin)
$ws.Range("D25").Value = "'2.237"
$ws.Range("E25").Value = "  -6.36%  "

# Row 26 (LidoDAOToken)
$ws.Range("D26").Value = "'2.529"
$ws.Range("E26").Value = "  -3.55%  "

# Row 27 (Monero)
$ws.Range("D27").Value = "'149.50"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").Value = "'18.88"
$ws.Range("E28").Value = "  -3.89%  "

# Row 29 (HuobiToken)
$ws.Range("D29").Value = "'4.836"
$ws.Range("E29").Value = "  -2.07%  "

# Row 30 (WrappedliquidstakedEther2.0)
$ws.Range("D30").Value = "1.729.75"
$ws.Range("E30").Value = "  -0.83%  "

# Row 31 (BitcoinCash)
$ws.Range("D31").Value = "'120.03"
$ws.Range("E31").Value = "  -4.20%  "

# Row 32 (ImmutableX)
$ws.Range("D32").Value = "'1.115"
$ws.Range("E32").Value = "  +1.75%  "

# Row 33 (Filecoin)
$ws.Range("D33").Value = "'5.929"
$ws.Range("E33").Value = "  -2.74%  "

# Row 34 (FraxShare)
$ws.Range("D34").Value = "'9.287"
$ws.Range("E34").Value = "  -5.74%  "

# Row 35 (Stellar)
$ws.Range("D35").Value = "'0.08105"
$ws.Range("E35").Value = "  -2.82%  "

# Row 36 (WEMIXTOKEN)
$ws.Range("D36").Value = "'1.604"
$ws.Range("E36").Value = "  -19.42%  "

# Row 37 / Row 38 swap: InternetComputer(DFINITY) <-> VeChain (identity swap with updated price/volume)
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02283"
$ws.Range("E37").Value = "  -7.04%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.215"
$ws.Range("E38").Value = "  -2.47%  "

# Row 39 (Hedera)
$ws.Range("D39").Value = "'0.06114"
$ws.Range("E39").Value = "  -4.39%  "

# Row 40 (Algorand)
$ws.Range("D40").Value = "'0.2111"
$ws.Range("E40").Value = "  -5.52%  "

# Row 41 (TrustWalletToken)
$ws.Range("D41").Value = "'1.215"
$ws.Range("E41").Value = "  -7.73%  "

# Row 42 (Aptos) - only volume changes
$ws.Range("E42").Value = "  -4.69%  "

# Row 43 (Frax) - only volume changes
$ws.Range("E43").Value = "  +0.01%  "

# Row 44 (TheSandbox)
$ws.Range("D44").Value = "'0.5937"
$ws.Range("E44").Value = "  -5.41%  "

# Row 45 (EnergySwap)
$ws.Range("D45").Value = "'13.62"
$ws.Range("E45").Value = "  -3.03%  "

# Row 46 (PancakeSwap)
$ws.Range("D46").Value = "'3.724"
$ws.Range("E46").Value = "  -1.26%  "

# Row 47 (Decentraland)
$ws.Range("D47").Value = "'0.5736"
$ws.Range("E47").Value = "  -5.89%  "

# Row 48 (NEARProtocol)
$ws.Range("D48").Value = "'1.940"
$ws.Range("E48").Value = "  -5.24%  "

# Row 49 (Quant)
$ws.Range("D49").Value = "'119.60"
$ws.Range("E49").Value = "  -4.40%  "

# Row 50 (EOS)
$ws.Range("D50").Value = "'1.155"
$ws.Range("E50").Value = "  -4.49%  "

# Row 51 (Cronos)
$ws.Range("D51").Value = "'0.06928"
$ws.Range("E51").Value = "  -3.81%  "
